$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New subject ids (column A)
$ws.Range("A19").Value = "SOC_XSIT_18"
$ws.Range("A20").Value = "SOC_XSIT_19"
$ws.Range("A21").Value = "SOC_XSIT_20"
$ws.Range("A22").Value = "SOC_XSIT_21"
$ws.Range("A23").Value = "SOC_XSIT_22"
$ws.Range("A24").Value = "SOC_XSIT_23"
$ws.Range("A25").Value = "SOC_XSIT_24"
$ws.Range("A26").Value = "SOC_XSIT_25"

# Run dates (column B)
$ws.Range("B19").Value = 41856
$ws.Range("B20").Value = 41856
$ws.Range("B21").Value = 41856
$ws.Range("B22").Value = 41856
$ws.Range("B23").Value = 41856
$ws.Range("B24").Value = 41857
$ws.Range("B25").Value = 41857
$ws.Range("B26").Value = 41858

# Gender (column E)
$ws.Range("E19").Value = "f"
$ws.Range("E20").Value = "m"
$ws.Range("E21").Value = "f"
$ws.Range("E22").Value = "f"
$ws.Range("E23").Value = "m"
$ws.Range("E24").Value = "f"
$ws.Range("E25").Value = "f"
$ws.Range("E26").Value = "f"

$ws.Range("G25").Select() | Out-Null
